$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.169.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.23%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.488.37'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.16%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.04%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.80%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.488.49'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.17%  '

$ws.Range('E10').Value = '  +0.84%  '

$ws.Range('E11').Value = '  +0.03%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.93'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.41%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.334'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.53%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.956.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.53%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.070.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.29%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000171'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.98%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.546.93'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.96%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.22%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.46%  '

$ws.Range('E23').Value = '  -0.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.53%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.67%  '

$ws.Range('E28').Value = '  +0.15%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.615.46'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0907'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '510.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.02%  '

$ws.Range('E32').Value = '  -3.85%  '

$ws.Range('E33').Value = '  -2.84%  '

$ws.Range('E34').Value = '  -4.04%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.18%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.117'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.68%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.96%  '

$ws.Range('E41').Value = '  -2.54%  '

$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('E43').Value = '  -2.42%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.59%  '

$ws.Range('E46').Value = '  -1.57%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.29%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.515'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.52%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.45'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.75%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0251'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.47%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0732'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.51%  '
